# Slide 18 ("Let's get started with NumPy") is reworked into a single
# big-text title slide with a URL, and the (empty) content placeholder
# below it is removed entirely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

$title = $s.Shapes.Item(1)

# Explicit position/size for the title placeholder (EMU 838200/2391952 /
# 10515600/1325563 expressed in points = EMU / 12700; nudged a hair above
# the exact boundary so float rounding lands on the same EMU value).
$title.Left   = 66.00003937007874
$title.Top    = 188.34271653543308
$title.Width  = 828.0000393700788
$title.Height = 104.37507874015748

# Replace the title text and make it large & bold.
$tr = $title.TextFrame.TextRange
$tr.Text = "https://tinyurl.com/sn-python"
$tr.Font.Size = 60
$tr.Font.Bold = $true

# The second shape ("Content Placeholder 2") was empty and is dropped.
$s.Shapes.Item(2).Delete()
